$d = $word.ActiveDocument

# Add template visibility: "{m:template myTemplate(...)}" -> "{m:template public myTemplate(...)}"
$d.Content.Find.Execute("myTemplate(a:notExisting", $true, $false, $false, $false, $false,
                         $true, 1, $false, "public myTemplate(a:notExisting", 2)

# Remove the stray _GoBack bookmark left over from the previous edit session
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}
